# Saldo.xlsx update — "Add files via upload"
#
# The source export was refreshed: a handful of account rows were added
# (new accounts, or accounts whose balance changed and got re-inserted at
# their new sorted position) and a few were removed (closed / zero'd out).
# All edits are on the single "Export" worksheet, which is a flat
# Conta/Nome/Saldo listing sorted by descending Saldo.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

function Insert-SaldoRow {
    param([int]$RowIndex, [string]$Conta, [string]$Nome, [double]$Saldo)
    $ws.Rows.Item($RowIndex).Insert()
    # Force column A to text so the leading zeros on the account number
    # survive, matching every other "Conta" cell in the sheet.
    $ws.Cells.Item($RowIndex, 1).Value = "'" + $Conta
    $ws.Cells.Item($RowIndex, 2).Value = $Nome
    $ws.Cells.Item($RowIndex, 3).Value = $Saldo
}

# --- Apply edits top-to-bottom, in original-sheet-row order ---------------

# New: 005749972 ALESSANDRA 389657.84 — inserted just above 005591536 (row 2)
Insert-SaldoRow 2 "005749972" "ALESSANDRA" 389657.84

# New: 004268684 PATRICIA 38441.84 — inserted just above 004211922 (now row 5)
Insert-SaldoRow 5 "004268684" "PATRICIA" 38441.84

# Removed: the old 004268684 PATRICIA 35000 row (now at row 7)
$ws.Rows.Item(7).Delete()

# Removed: 004452476 IVONE 6309.86 (now at row 9)
$ws.Rows.Item(9).Delete()

# New: 004482102 NATALIA 4000 — inserted just below 004397124 (row 9), i.e. at row 10
Insert-SaldoRow 10 "004482102" "NATALIA" 4000

# New: 004448303 NASSIM 2147.39 — inserted right after the NATALIA row, at row 11
Insert-SaldoRow 11 "004448303" "NASSIM" 2147.39

# New: 005186167 ANDREA 97.19 — inserted just above 004211911 (now row 29)
Insert-SaldoRow 29 "005186167" "ANDREA" 97.19

# Removed: 004207278 CESAR 11.92 (now at row 199)
$ws.Rows.Item(199).Delete()
